# Update test cases for data-brain-test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getConceptModelDataByCondition")

# Update the description column (B) for rows 2 and 3 to the new text
$ws.Range("B2").Value = "good request, data retrieved(no schema check)"
$ws.Range("B3").Value = "good request, data retrieved(no schema check)"

# Widen column B to fit the longer text (stored width ends up as 38 after
# Excel's internal padding is re-added on save; ColumnWidth itself reports
# the value before that padding, so back it out here).
$ws.Columns.Item(2).ColumnWidth = 37 + 1/6

# Move the active selection to F9
$ws.Range("F9").Select()
